$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3, 1, 31.29063333333333, 93.8719, 0.02026792284095206, 0.02026792284095206, 3, 1, 59.02895333333333, 177.08686, 0.2011980443121526, 0.2011980443121526, 1847.053334803778, 16623.480013234, 0.004077866437869162, 0.004077866437869162),
    @(3, 3, 1, 31.29063333333333, 93.8719, 0.02026792284095206, 0.02026792284095206, 3, 1, 71.14312966666667, 213.429389, 0.2424887745230654, 0.2424887745230655, 2226.113584585456, 20035.0222612691, 0.004914743771830511, 0.004914743771830512),
    @(4, 3, 1, 31.29063333333333, 93.8719, 0.02026792284095206, 0.02026792284095206, 3, 1, 68.402828, 205.208484, 0.233148555782522, 0.233148555782522, 2140.367809911067, 19263.3102891996, 0.004725436939079562, 0.004725436939079562),
    @(5, 3, 1, 31.29063333333333, 93.8719, 0.02026792284095206, 0.02026792284095206, 3, 1, 74.34717300000001, 223.041519, 0.2534096398976854, 0.2534096398976854, 2326.370129712901, 20937.3311674161, 0.005136087028599734, 0.005136087028599734),
    @(6, 3, 1, 31.29063333333333, 93.8719, 0.02026792284095206, 0.02026792284095206, 3, 1, 20.46522766666667, 61.395683, 0.0697549854845745, 0.0697549854845745, 640.3699350008555, 5763.3294150077, 0.001413788663573086, 0.001413788663573086),
    @(7, 3, 1, 170.232249, 510.696747, 0.1102647572204378, 0.1102647572204378, 3, 1, 59.02895333333333, 177.08686, 0.2011980443121526, 0.2011980443121526, 10048.63148204938, 90437.68333844442, 0.02218505350930639, 0.02218505350930639),
    @(8, 3, 1, 170.232249, 510.696747, 0.1102647572204378, 0.1102647572204378, 3, 1, 71.14312966666667, 213.429389, 0.2424887745230654, 0.2424887745230655, 12110.85496405529, 108997.6946764976, 0.02673796585146729, 0.02673796585146729),
    @(9, 3, 1, 170.232249, 510.696747, 0.1102647572204378, 0.1102647572204378, 3, 1, 68.402828, 205.208484, 0.233148555782522, 0.233148555782522, 11644.36724840017, 104799.3052356015, 0.02570806889965548, 0.02570806889965548),
    @(10, 3, 1, 170.232249, 510.696747, 0.1102647572204378, 0.1102647572204378, 3, 1, 74.34717300000001, 223.041519, 0.2534096398976854, 0.2534096398976854, 12656.28646658208, 113906.5781992387, 0.02794215242063685, 0.02794215242063685),
    @(11, 3, 1, 170.232249, 510.696747, 0.1102647572204378, 0.1102647572204378, 3, 1, 20.46522766666667, 61.395683, 0.0697549854845745, 0.0697549854845745, 3483.841731993689, 31354.5755879432, 0.007691516539371768, 0.007691516539371768),
    @(12, 3, 1, 603.9765116666666, 1811.929535, 0.391214495590503, 0.391214495590503, 3, 1, 59.02895333333333, 177.08686, 0.2011980443121526, 0.2011980443121526, 35652.10132160112, 320868.9118944101, 0.07871159141937446, 0.07871159141937446),
    @(13, 3, 1, 603.9765116666666, 1811.929535, 0.391214495590503, 0.391214495590503, 3, 1, 71.14312966666667, 213.429389, 0.2424887745230654, 0.2424887745230655, 42968.77928512268, 386719.0135661041, 0.09486512361140026, 0.09486512361140027),
    @(14, 3, 1, 603.9765116666666, 1811.929535, 0.391214495590503, 0.391214495590503, 3, 1, 68.402828, 205.208484, 0.233148555782522, 0.233148555782522, 41313.70144357499, 371823.3129921749, 0.09121109464811358, 0.09121109464811358),
    @(15, 3, 1, 603.9765116666666, 1811.929535, 0.391214495590503, 0.391214495590503, 3, 1, 74.34717300000001, 223.041519, 0.2534096398976854, 0.2534096398976854, 44903.94620081819, 404135.5158073637, 0.099137524450344, 0.099137524450344),
    @(16, 3, 1, 603.9765116666666, 1811.929535, 0.391214495590503, 0.391214495590503, 3, 1, 20.46522766666667, 61.395683, 0.0697549854845745, 0.0697549854845745, 12360.51681657749, 111244.6513491974, 0.02728916146127067, 0.02728916146127067),
    @(17, 3, 1, 694.886846, 2084.660538, 0.4500999653118942, 0.4500999653118942, 3, 1, 59.02895333333333, 177.08686, 0.2011980443121526, 0.2011980443121526, 41018.44320448118, 369165.9888403307, 0.09055923276572085, 0.09055923276572085),
    @(18, 3, 1, 694.886846, 2084.660538, 0.4500999653118942, 0.4500999653118942, 3, 1, 71.14312966666667, 213.429389, 0.2424887745230654, 0.2424887745230655, 49436.42498863903, 444927.8248977513, 0.1091441890013555, 0.1091441890013555),
    @(19, 3, 1, 694.886846, 2084.660538, 0.4500999653118942, 0.4500999653118942, 3, 1, 68.402828, 205.208484, 0.233148555782522, 0.233148555782522, 47532.22540640049, 427790.0286576044, 0.1049401568702314, 0.1049401568702314),
    @(20, 3, 1, 694.886846, 2084.660538, 0.4500999653118942, 0.4500999653118942, 3, 1, 74.34717300000001, 223.041519, 0.2534096398976854, 0.2534096398976854, 51662.87255498637, 464965.8529948773, 0.1140596701276478, 0.1140596701276478),
    @(21, 3, 1, 694.886846, 2084.660538, 0.4500999653118942, 0.4500999653118942, 3, 1, 20.46522766666667, 61.395683, 0.0697549854845745, 0.0697549854845745, 14221.01750596194, 127989.1575536575, 0.03139671654693867, 0.03139671654693867),
    @(22, 3, 1, 43.46379233333334, 130.391377, 0.02815285903621308, 0.02815285903621308, 3, 1, 59.02895333333333, 177.08686, 0.2011980443121526, 0.2011980443121526, 2565.622169334025, 23090.59952400622, 0.005664300179881786, 0.005664300179881786),
    @(23, 3, 1, 43.46379233333334, 130.391377, 0.02815285903621308, 0.02815285903621308, 3, 1, 71.14312966666667, 213.429389, 0.2424887745230654, 0.2424887745230655, 3092.150213775406, 27829.35192397865, 0.00682675228701192, 0.00682675228701192),
    @(24, 3, 1, 43.46379233333334, 130.391377, 0.02815285903621308, 0.02815285903621308, 3, 1, 68.402828, 205.208484, 0.233148555782522, 0.233148555782522, 2973.046311204719, 26757.41680084247, 0.006563798425442005, 0.006563798425442004),
    @(25, 3, 1, 43.46379233333334, 130.391377, 0.02815285903621308, 0.02815285903621308, 3, 1, 74.34717300000001, 223.041519, 0.2534096398976854, 0.2534096398976854, 3231.410087842408, 29082.69079058167, 0.007134205870457057, 0.007134205870457056),
    @(26, 3, 1, 43.46379233333334, 130.391377, 0.02815285903621308, 0.02815285903621308, 3, 1, 20.46522766666667, 61.395683, 0.0697549854845745, 0.0697549854845745, 889.496405358388, 8005.467648225491, 0.001963802273420316, 0.001963802273420315)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 5).Value = $row[1]
    $ws.Cells.Item($r, 6).Value = $row[2]
    $ws.Cells.Item($r, 7).Value = $row[3]
    $ws.Cells.Item($r, 8).Value = $row[4]
    $ws.Cells.Item($r, 9).Value = $row[5]
    $ws.Cells.Item($r, 10).Value = $row[6]
    $ws.Cells.Item($r, 11).Value = $row[7]
    $ws.Cells.Item($r, 12).Value = $row[8]
    $ws.Cells.Item($r, 13).Value = $row[9]
    $ws.Cells.Item($r, 14).Value = $row[10]
    $ws.Cells.Item($r, 15).Value = $row[11]
    $ws.Cells.Item($r, 16).Value = $row[12]
    $ws.Cells.Item($r, 17).Value = $row[13]
    $ws.Cells.Item($r, 18).Value = $row[14]
    $ws.Cells.Item($r, 19).Value = $row[15]
    $ws.Cells.Item($r, 20).Value = $row[16]
}
